# Replace the body text of the response paragraph with the revised wording.
$d = $word.ActiveDocument

$oldText = "In my current and previous position I never worked directly with passengers data, so I am curious on how people working with people’s data handle the situation to ensure they are not easily identifiable. The only example were I had to mask some information to protect the identity of the client (customer) was when presenting data about a project or product to a potential new client for the GE Shop Services for turbojet/turbofan engines. We needed to show the data/trends and tools we were able to offer to them as well as some of the advantages our current customers experience such as: improvements on engine performance & improvements on fuel consumption e, all this by making sure that the data coming from each of the customers were not identifiable (labels, locations and/or anything that could potentially lead to identify the customers/airlines. Some of the processes we followed were: change the label/names and location of the engines and instead plot/present everything in terms of ambient conditions, if possible take a random sample from all the data points available removing any labels."
$newText = "Unfortunately, I have not worked directly with labels directly related to people, so I am curious on how people’s data is handled to ensure they are not easily identifiable. The only example that I have when I had to mask (somehow) some information to protect the identity of the client/customer was when presenting data about a project or product to a potential new client for the GE Shop Services for turbojet/turbofan engines. We needed to show the data/trends and tools we were able to offer to them as well as some of the advantages our current customers experience such as: improvements on engine performance & improvements on fuel consumption. All this by making sure that the data coming from each of the customers was not identifiable (labels, locations and/or anything that could potentially lead to identify the customers/airlines). This are some of the process that I remember (kind of) we followed: Change the label/names and location of the engines and instead plot/present everything in terms of ambient conditions, in other cases we would remove all together name labels and/or location labels and if  possible take a random sample from all the data points available removing any labels."

$found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

if (-not $found) {
    throw "Could not locate the paragraph text to replace."
}

# Remove the leftover "_GoBack" last-edit-position bookmark, matching the
# author's resave of the document after the text was updated.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$d.Save()
